$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1447").NumberFormat = "@"
$ws.Range("D1447").Value = "38"
$ws.Range("E1447").NumberFormat = "@"
$ws.Range("E1447").Value = "6"
$ws.Range("F1447").NumberFormat = "@"
$ws.Range("F1447").Value = "29"
$ws.Range("H1447").NumberFormat = "@"
$ws.Range("H1447").Value = "3"
$ws.Range("B1448").NumberFormat = "@"
$ws.Range("B1448").Value = "2"
$ws.Range("C1448").Value = "Miami Dolphins"
$ws.Range("D1448").NumberFormat = "@"
$ws.Range("D1448").Value = "35"
$ws.Range("E1448").NumberFormat = "@"
$ws.Range("E1448").Value = "8"
$ws.Range("F1448").NumberFormat = "@"
$ws.Range("F1448").Value = "23"
$ws.Range("G1448").NumberFormat = "@"
$ws.Range("G1448").Value = "1"
$ws.Range("B1449").NumberFormat = "@"
$ws.Range("B1449").Value = "2"
$ws.Range("C1449").Value = "Tampa Bay Buccaneers"
$ws.Range("D1449").NumberFormat = "@"
$ws.Range("D1449").Value = "35"
$ws.Range("E1449").NumberFormat = "@"
$ws.Range("E1449").Value = "6"
$ws.Range("F1449").NumberFormat = "@"
$ws.Range("F1449").Value = "25"
$ws.Range("G1449").NumberFormat = "@"
$ws.Range("G1449").Value = "0"
$ws.Range("H1449").NumberFormat = "@"
$ws.Range("H1449").Value = "4"
$ws.Range("B1450").NumberFormat = "@"
$ws.Range("B1450").Value = "4"
$ws.Range("D1450").NumberFormat = "@"
$ws.Range("D1450").Value = "34"
$ws.Range("E1450").NumberFormat = "@"
$ws.Range("E1450").Value = "12"
$ws.Range("F1450").NumberFormat = "@"
$ws.Range("F1450").Value = "17"
$ws.Range("H1450").NumberFormat = "@"
$ws.Range("H1450").Value = "5"
$ws.Range("C1451").Value = "Detroit Lions"
$ws.Range("D1451").NumberFormat = "@"
$ws.Range("D1451").Value = "33"
$ws.Range("E1451").NumberFormat = "@"
$ws.Range("E1451").Value = "10"
$ws.Range("F1451").NumberFormat = "@"
$ws.Range("F1451").Value = "22"
$ws.Range("H1451").NumberFormat = "@"
$ws.Range("H1451").Value = "1"
$ws.Range("C1452").Value = "Cincinnati Bengals"
$ws.Range("D1452").NumberFormat = "@"
$ws.Range("D1452").Value = "31"
$ws.Range("E1452").NumberFormat = "@"
$ws.Range("E1452").Value = "13"
$ws.Range("F1452").NumberFormat = "@"
$ws.Range("F1452").Value = "15"
$ws.Range("H1452").NumberFormat = "@"
$ws.Range("H1452").Value = "3"
$ws.Range("C1453").Value = "Oakland Raiders"
$ws.Range("D1453").NumberFormat = "@"
$ws.Range("D1453").Value = "31"
$ws.Range("E1453").NumberFormat = "@"
$ws.Range("E1453").Value = "7"
$ws.Range("F1453").NumberFormat = "@"
$ws.Range("F1453").Value = "24"
$ws.Range("G1453").NumberFormat = "@"
$ws.Range("G1453").Value = "0"
$ws.Range("H1453").NumberFormat = "@"
$ws.Range("H1453").Value = "0"
$ws.Range("C1454").Value = "Atlanta Falcons"
$ws.Range("D1454").NumberFormat = "@"
$ws.Range("D1454").Value = "30"
$ws.Range("F1454").NumberFormat = "@"
$ws.Range("F1454").Value = "19"
$ws.Range("H1454").NumberFormat = "@"
$ws.Range("H1454").Value = "2"
$ws.Range("C1455").Value = "Carolina Panthers"
$ws.Range("D1455").NumberFormat = "@"
$ws.Range("D1455").Value = "30"
$ws.Range("E1455").NumberFormat = "@"
$ws.Range("E1455").Value = "18"
$ws.Range("F1455").NumberFormat = "@"
$ws.Range("F1455").Value = "11"
$ws.Range("B1456").NumberFormat = "@"
$ws.Range("B1456").Value = "8"
$ws.Range("C1456").Value = "Kansas City Chiefs"
$ws.Range("D1456").NumberFormat = "@"
$ws.Range("D1456").Value = "30"
$ws.Range("E1456").NumberFormat = "@"
$ws.Range("E1456").Value = "12"
$ws.Range("F1456").NumberFormat = "@"
$ws.Range("F1456").Value = "17"
$ws.Range("H1456").NumberFormat = "@"
$ws.Range("H1456").Value = "1"
$ws.Range("B1457").NumberFormat = "@"
$ws.Range("B1457").Value = "8"
$ws.Range("C1457").Value = "New York Jets"
$ws.Range("D1457").NumberFormat = "@"
$ws.Range("D1457").Value = "30"
$ws.Range("E1457").NumberFormat = "@"
$ws.Range("E1457").Value = "10"
$ws.Range("F1457").NumberFormat = "@"
$ws.Range("F1457").Value = "18"
$ws.Range("H1457").NumberFormat = "@"
$ws.Range("H1457").Value = "2"
$ws.Range("B1458").NumberFormat = "@"
$ws.Range("B1458").Value = "8"
$ws.Range("C1458").Value = "Washington Redskins"
$ws.Range("D1458").NumberFormat = "@"
$ws.Range("D1458").Value = "30"
$ws.Range("E1458").NumberFormat = "@"
$ws.Range("E1458").Value = "9"
$ws.Range("F1458").NumberFormat = "@"
$ws.Range("F1458").Value = "19"
$ws.Range("B1459").NumberFormat = "@"
$ws.Range("B1459").Value = "13"
$ws.Range("C1459").Value = "Seattle Seahawks"
$ws.Range("D1459").NumberFormat = "@"
$ws.Range("D1459").Value = "29"
$ws.Range("E1459").NumberFormat = "@"
$ws.Range("E1459").Value = "12"
$ws.Range("F1459").NumberFormat = "@"
$ws.Range("F1459").Value = "12"
$ws.Range("G1459").NumberFormat = "@"
$ws.Range("G1459").Value = "1"
$ws.Range("H1459").NumberFormat = "@"
$ws.Range("H1459").Value = "4"
$ws.Range("B1460").NumberFormat = "@"
$ws.Range("B1460").Value = "14"
$ws.Range("C1460").Value = "Houston Texans"
$ws.Range("D1460").NumberFormat = "@"
$ws.Range("D1460").Value = "28"
$ws.Range("E1460").NumberFormat = "@"
$ws.Range("E1460").Value = "6"
$ws.Range("F1460").NumberFormat = "@"
$ws.Range("F1460").Value = "22"
$ws.Range("H1460").NumberFormat = "@"
$ws.Range("H1460").Value = "0"
$ws.Range("C1461").Value = "Cleveland Browns"
$ws.Range("D1461").NumberFormat = "@"
$ws.Range("D1461").Value = "27"
$ws.Range("E1461").NumberFormat = "@"
$ws.Range("E1461").Value = "8"
$ws.Range("F1461").NumberFormat = "@"
$ws.Range("F1461").Value = "17"
$ws.Range("H1461").NumberFormat = "@"
$ws.Range("H1461").Value = "2"
$ws.Range("C1462").Value = "Philadelphia Eagles"
$ws.Range("D1462").NumberFormat = "@"
$ws.Range("D1462").Value = "27"
$ws.Range("E1462").NumberFormat = "@"
$ws.Range("E1462").Value = "9"
$ws.Range("F1462").NumberFormat = "@"
$ws.Range("F1462").Value = "17"
$ws.Range("G1462").NumberFormat = "@"
$ws.Range("G1462").Value = "1"
$ws.Range("H1462").NumberFormat = "@"
$ws.Range("H1462").Value = "0"
$ws.Range("C1463").Value = "Jacksonville Jaguars"
$ws.Range("D1463").NumberFormat = "@"
$ws.Range("D1463").Value = "26"
$ws.Range("E1463").NumberFormat = "@"
$ws.Range("E1463").Value = "11"
$ws.Range("F1463").NumberFormat = "@"
$ws.Range("F1463").Value = "14"
$ws.Range("H1463").NumberFormat = "@"
$ws.Range("H1463").Value = "1"
$ws.Range("C1464").Value = "Los Angeles Chargers"
$ws.Range("D1464").NumberFormat = "@"
$ws.Range("D1464").Value = "26"
$ws.Range("E1464").NumberFormat = "@"
$ws.Range("E1464").Value = "9"
$ws.Range("F1464").NumberFormat = "@"
$ws.Range("F1464").Value = "15"
$ws.Range("H1464").NumberFormat = "@"
$ws.Range("H1464").Value = "2"
$ws.Range("C1465").Value = "Indianapolis Colts"
$ws.Range("D1465").NumberFormat = "@"
$ws.Range("D1465").Value = "24"
$ws.Range("F1465").NumberFormat = "@"
$ws.Range("F1465").Value = "16"
$ws.Range("H1465").NumberFormat = "@"
$ws.Range("H1465").Value = "2"
$ws.Range("C1466").Value = "Green Bay Packers"
$ws.Range("D1466").NumberFormat = "@"
$ws.Range("D1466").Value = "23"
$ws.Range("E1466").NumberFormat = "@"
$ws.Range("E1466").Value = "11"
$ws.Range("H1466").NumberFormat = "@"
$ws.Range("H1466").Value = "0"
$ws.Range("C1467").Value = "Pittsburgh Steelers"
$ws.Range("D1467").NumberFormat = "@"
$ws.Range("D1467").Value = "23"
$ws.Range("E1467").NumberFormat = "@"
$ws.Range("E1467").Value = "5"
$ws.Range("F1467").NumberFormat = "@"
$ws.Range("F1467").Value = "17"
$ws.Range("B1468").NumberFormat = "@"
$ws.Range("B1468").Value = "22"
$ws.Range("C1468").Value = "Los Angeles Rams"
$ws.Range("D1468").NumberFormat = "@"
$ws.Range("D1468").Value = "22"
$ws.Range("E1468").NumberFormat = "@"
$ws.Range("E1468").Value = "7"
$ws.Range("F1468").NumberFormat = "@"
$ws.Range("F1468").Value = "13"
$ws.Range("H1468").NumberFormat = "@"
$ws.Range("H1468").Value = "2"
$ws.Range("B1469").NumberFormat = "@"
$ws.Range("B1469").Value = "22"
$ws.Range("C1469").Value = "Minnesota Vikings"
$ws.Range("D1469").NumberFormat = "@"
$ws.Range("D1469").Value = "22"
$ws.Range("E1469").NumberFormat = "@"
$ws.Range("E1469").Value = "3"
$ws.Range("F1469").NumberFormat = "@"
$ws.Range("F1469").Value = "19"
$ws.Range("B1470").NumberFormat = "@"
$ws.Range("B1470").Value = "22"
$ws.Range("C1470").Value = "New Orleans Saints"
$ws.Range("D1470").NumberFormat = "@"
$ws.Range("D1470").Value = "22"
$ws.Range("E1470").NumberFormat = "@"
$ws.Range("E1470").Value = "7"
$ws.Range("F1470").NumberFormat = "@"
$ws.Range("F1470").Value = "14"
$ws.Range("G1470").NumberFormat = "@"
$ws.Range("G1470").Value = "0"
$ws.Range("B1471").NumberFormat = "@"
$ws.Range("B1471").Value = "22"
$ws.Range("C1471").Value = "Tennessee Titans"
$ws.Range("D1471").NumberFormat = "@"
$ws.Range("D1471").Value = "22"
$ws.Range("E1471").NumberFormat = "@"
$ws.Range("E1471").Value = "5"
$ws.Range("F1471").NumberFormat = "@"
$ws.Range("F1471").Value = "17"
$ws.Range("G1471").NumberFormat = "@"
$ws.Range("G1471").Value = "0"
$ws.Range("C1472").Value = "Buffalo Bills"
$ws.Range("D1472").NumberFormat = "@"
$ws.Range("D1472").Value = "21"
$ws.Range("E1472").NumberFormat = "@"
$ws.Range("E1472").Value = "11"
$ws.Range("F1472").NumberFormat = "@"
$ws.Range("F1472").Value = "7"
$ws.Range("G1472").NumberFormat = "@"
$ws.Range("G1472").Value = "1"
$ws.Range("H1472").NumberFormat = "@"
$ws.Range("H1472").Value = "2"
$ws.Range("B1473").NumberFormat = "@"
$ws.Range("B1473").Value = "26"
$ws.Range("C1473").Value = "Dallas Cowboys"
$ws.Range("D1473").NumberFormat = "@"
$ws.Range("D1473").Value = "21"
$ws.Range("E1473").NumberFormat = "@"
$ws.Range("E1473").Value = "10"
$ws.Range("F1473").NumberFormat = "@"
$ws.Range("F1473").Value = "11"
$ws.Range("H1473").NumberFormat = "@"
$ws.Range("H1473").Value = "0"
$ws.Range("B1474").NumberFormat = "@"
$ws.Range("B1474").Value = "28"
$ws.Range("C1474").Value = "Baltimore Ravens"
$ws.Range("D1474").NumberFormat = "@"
$ws.Range("D1474").Value = "20"
$ws.Range("E1474").NumberFormat = "@"
$ws.Range("E1474").Value = "10"
$ws.Range("F1474").NumberFormat = "@"
$ws.Range("F1474").Value = "9"
$ws.Range("G1474").NumberFormat = "@"
$ws.Range("G1474").Value = "1"
$ws.Range("C1475").Value = "Chicago Bears"
$ws.Range("D1475").NumberFormat = "@"
$ws.Range("D1475").Value = "19"
$ws.Range("E1475").NumberFormat = "@"
$ws.Range("E1475").Value = "10"
$ws.Range("F1475").NumberFormat = "@"
$ws.Range("F1475").Value = "9"
$ws.Range("G1475").NumberFormat = "@"
$ws.Range("G1475").Value = "0"
$ws.Range("H1475").NumberFormat = "@"
$ws.Range("H1475").Value = "0"
$ws.Range("B1476").NumberFormat = "@"
$ws.Range("B1476").Value = "29"
$ws.Range("C1476").Value = "Denver Broncos"
$ws.Range("D1476").NumberFormat = "@"
$ws.Range("D1476").Value = "19"
$ws.Range("E1476").NumberFormat = "@"
$ws.Range("E1476").Value = "6"
$ws.Range("F1476").NumberFormat = "@"
$ws.Range("F1476").Value = "10"
$ws.Range("G1476").NumberFormat = "@"
$ws.Range("G1476").Value = "1"
$ws.Range("H1476").NumberFormat = "@"
$ws.Range("H1476").Value = "1"
$ws.Range("C1477").Value = "San Francisco 49ers"
$ws.Range("D1477").NumberFormat = "@"
$ws.Range("D1477").Value = "17"
$ws.Range("F1477").NumberFormat = "@"
$ws.Range("F1477").Value = "10"
$ws.Range("H1477").NumberFormat = "@"
$ws.Range("H1477").Value = "2"
$ws.Range("B1478").NumberFormat = "@"
$ws.Range("B1478").Value = "31"
$ws.Range("C1478").Value = "New England Patriots"
$ws.Range("D1478").NumberFormat = "@"
$ws.Range("D1478").Value = "13"
$ws.Range("E1478").NumberFormat = "@"
$ws.Range("E1478").Value = "5"
$ws.Range("F1478").NumberFormat = "@"
$ws.Range("F1478").Value = "4"
$ws.Range("H1478").NumberFormat = "@"
$ws.Range("H1478").Value = "4"

Write-Host "Applied 173 cell updates"
